# overworld 1 lesson and tutorial
# Adds a new "latitude"/"Latitude" row in the middle of the language table
# (inserted as new row 25, pushing the season/atmosphere/climate/region/
# weather/intro rows down by one), then appends 19 new rows describing the
# "overworld_1" tutorial + sun-illustration strings at the bottom of the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "latitude" row at position 25 -------------------------
# (this shifts every row from the old 25 onward down by one, exactly like
# the diff shows for season_title..intro_4_2)
$ws.Rows(25).Insert()
$ws.Range("A25").Value = "latitude"
$ws.Range("B25").Value = "Latitude"

# --- Append the new overworld_1 tutorial rows at the bottom ---------------
# Old last row was 81 (intro_4_2); after the insert above it is row 82, so
# the new content starts at row 83.
$newRows = @(
    @("sunIllustrate_title",        "Sunlight Direction"),
    @("sunIllustrate_hot",          "Hot"),
    @("sunIllustrate_warm",         "Warm"),
    @("overworld_1_intro_0",        "This is the map of Earth. The frogs are very particular with their choice of habitat."),
    @("overworld_1_intro_1",        "We’ll need to determine where to place the frogs by looking at their criteria."),
    @("overworld_1_criteria_0",     "In this case, we need to find a place where it’s fairly warm and humid."),
    @("overworld_1_hud_0",          "On the lefthand side of the map are the latitude values. This is the angular distance relative to the earth’s equator (middle of the map)."),
    @("overworld_1_hud_1",          "Now let’s view the temperature readings of Earth."),
    @("overworld_1_temp_0",         "As you can see, the temperature is consistently hot starting from the equator, and gets colder further north or south."),
    @("overworld_1_humid_0",        "Next is the humidity readings of Earth. The percentage tells us how much water vapor is present on air."),
    @("overworld_1_humid_1",        "Notice how humidity tends to be higher in large forest and jungle areas, such as the Amazon rainforest."),
    @("overworld_1_hotspot_0",      "Now go ahead and find a place for the frogs to land. Simply click around the map to find the spot."),
    @("overworld_1_hotspot_1",      "Remember to make use of the temperature and humidity reading. "),
    @("overworld_1_hotspot_2",      "Hint: we are looking for a temperate climate, perhaps somewhere in North America..."),
    @("overworld_1_analyze_0",      "Looks like the temperature is too low. We will need to change the time of the year to land on this spot."),
    @("overworld_1_analyze_1",      "Since the earth rotate at a slightly tilted axis around the sun, the atmosphere can change throughout the year."),
    @("overworld_1_analyze_2",      "Why don’t we adjust the time by selecting a different season."),
    @("overworld_1_investigate_0",  "Now you just need to pick a particular location on the land for the frogs."),
    @("overworld_1_investigate_1",  "When the majority of the frogs approve, we can finally launch the expedition!")
)

$row = 83
foreach ($pair in $newRows) {
    $ws.Range("A$row").Value = $pair[0]
    $ws.Range("B$row").Value = $pair[1]
    $row = $row + 1
}

# Row 86 (overworld_1_intro_0's value) is vertically centered, matching the
# style already used on the intro_4_1 / intro_4_2 rows above it.
$ws.Range("B86").VerticalAlignment = -4108

# --- Selection / scroll state, matching the saved view in the diff --------
$null = $ws.Range("B93").Select()
